# Refresh the crypto price snapshot (Price / Volume(1h) columns) plus the
# BKEXToken <-> KickToken rank swap in rows 41-42, matching the
# "Updated symbol list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new literal-text value.
$updates = @(
    @('D2', '255.36'),
    @('E2', '3.93%'),
    @('D3', '28.20'),
    @('E3', '-4.41%'),
    @('D4', '5.207'),
    @('E4', '-1.83%'),
    @('E5', '2.19%'),
    @('D6', '6.698'),
    @('E6', '0.75%'),
    @('E7', '1.49%'),
    @('D8', '0.9582'),
    @('E8', '12.37%'),
    @('D9', '0.1410'),
    @('E9', '2.26%'),
    @('D10', '0.07155'),
    @('E10', '0.77%'),
    @('D11', '0.03209'),
    @('E11', '0.23%'),
    @('D12', '0.09221'),
    @('E12', '-1.30%'),
    @('D13', '0.001538'),
    @('E13', '-0.23%'),
    @('E14', '-94.01%'),
    @('D15', '0.005859'),
    @('E15', '-3.76%'),
    @('D16', '3.500'),
    @('E16', '-0.41%'),
    @('D17', '3.211'),
    @('E17', '0.22%'),
    @('E18', '1.39%'),
    @('E19', '0.50%'),
    @('D20', '0.03475'),
    @('E20', '3.80%'),
    @('E21', '0.32%'),
    @('D22', '3.527'),
    @('E22', '0.52%'),
    @('D23', '0.04178'),
    @('E23', '0.91%'),
    @('E24', '-0.93%'),
    @('D25', '0.001226'),
    @('E25', '0.42%'),
    @('D26', '0.004557'),
    @('E26', '9.49%'),
    @('E27', '-0.04%'),
    @('E28', '1.13%'),
    @('D40', '0.03817'),
    @('E40', '1.30%'),
    @('B41', 'KickToken'),
    @('C41', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'),
    @('D41', '0.005619'),
    @('E41', '-2.32%'),
    @('B42', 'BKEXToken'),
    @('C42', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'),
    @('D42', '0.1100'),
    @('E42', '2.93%'),
    @('D43', '0.002300'),
    @('E43', '-4.88%'),
    @('D44', '0.009733'),
    @('E44', '5.80%'),
    @('D45', '0.00005387'),
    @('E45', '1.79%'),
    @('E46', '-0.07%'),
    @('D47', '0.08999'),
    @('E47', '11.19%'),
    @('E48', '-3.41%'),
    @('E49', '-0.07%'),
    @('E50', '-0.07%')
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $cell = $ws.Range($cellRef)
    # Leading apostrophe forces literal-text storage so numeric-looking
    # price/percentage strings are not reinterpreted as numbers.
    $cell.Value = "'" + $newValue
    # Re-apply the default style so the apostrophe-prefix quoting does
    # not leave a visible "quote prefix" cell-format behind.
    $cell.Style = "Normal"
}

